$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row above the old row 4 ("Number of disability persons")
#    so we end up with:
#      row1 = title (merged)
#      row2 = "(End of year, persons)"
#      row3 = year headers
#      row4 = NEW "family with disabilities Persons " data
#      row5 = old row4 data, relabeled "disabilities Persons "
#      row6 = old row5, source note
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# Copy formatting from row 5 (the data row that used to be row 4) into the
# freshly inserted row 4 so fonts/fills/number formats match without having
# to rebuild every style from scratch.
$ws.Range("A5:I5").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Row 1 - title, merged across A1:I1
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").UnMerge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Mestia Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51
$ws.Range("A1:I1").Select()

# ---------------------------------------------------------------------------
# 3. Row 2 - unchanged text, just restore default row height
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Row 3 - A3 font becomes Sylfaen 11 (still themed black, top border)
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.ThemeColor = 1
$ws.Rows.Item(3).RowHeight = 18.75

# ---------------------------------------------------------------------------
# 5. Row 4 (new) - "family with disabilities Persons "
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value2 = 207
$ws.Range("C4").Value2 = 204
$ws.Range("D4").Value2 = 208
$ws.Range("E4").Value2 = 204
$ws.Range("F4").Value2 = 209
$ws.Range("G4").Value2 = 210
$ws.Range("H4").Value2 = 211
$ws.Range("I4").Value2 = 212
$ws.Range("B4:I4").NumberFormat = "#\ ##0"
$ws.Range("B4:I4").HorizontalAlignment = -4142
$ws.Range("A4:I4").Borders.Item(9).LineStyle = -4142
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 6. Row 5 (previously row 4) - "disabilities Persons ", new values
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value2 = 239
$ws.Range("C5").Value2 = 232
$ws.Range("D5").Value2 = 238
$ws.Range("E5").Value2 = 234
$ws.Range("F5").Value2 = 240
$ws.Range("G5").Value2 = 240
$ws.Range("H5").Value2 = 235
$ws.Range("I5").Value2 = 235
$ws.Range("B5:I5").NumberFormat = "#\ ##0"
$ws.Range("B5:H5").HorizontalAlignment = -4142
$ws.Range("A5:I5").Borders.Item(8).LineStyle = -4142
$ws.Range("I5").Borders.Item(9).LineStyle = -4142
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------------
# 7. Row 6 (previously row 5) - source note, merged A6:H6
# ---------------------------------------------------------------------------
$ws.Range("A6:H6").Borders.Item(8).LineStyle = -4142
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Sheet-wide cosmetics
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.83
